$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '34.532.33'
$ws.Range("E2").Value = '  +14.03%  '
$ws.Range("D3").Value = '1.801.79'
$ws.Range("E3").Value = '  +7.54%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.996'
$ws.Range("E4").Value = '  -0.16%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '233.06'
$ws.Range("E5").Value = '  +5.99%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.551'
$ws.Range("E6").Value = '  +5.26%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.997'
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '31.65'
$ws.Range("E8").Value = '  +6.20%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '46.28'
$ws.Range("E9").Value = '  +5.34%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.283'
$ws.Range("E10").Value = '  +6.73%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0679'
$ws.Range("E11").Value = '  +9.50%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0926'
$ws.Range("E12").Value = '  +2.04%  '
$ws.Range("D13").Value = '2.057.87'
$ws.Range("E13").Value = '  +7.46%  '
$ws.Range("D14").Value = '1.788.64'
$ws.Range("E14").Value = '  +7.20%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.642'
$ws.Range("E15").Value = '  +3.29%  '
$ws.Range("D16").Value = '34.491.64'
$ws.Range("E16").Value = '  +13.88%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '10.24'
$ws.Range("E17").Value = '  -5.17%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '4.34'
$ws.Range("E18").Value = '  +8.21%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '71.03'
$ws.Range("E19").Value = '  +8.00%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '264.29'
$ws.Range("E20").Value = '  +6.56%  '
$ws.Range("D21").Value = '0.0₃0761'
$ws.Range("E21").Value = '  +5.86%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.995'
$ws.Range("E22").Value = '  -0.34%  '
$ws.Range("B23").Value = 'Uniswap'
$ws.Range("C23").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.41'
$ws.Range("E23").Value = '  +1.89%  '
$ws.Range("B24").Value = 'Avalanche'
$ws.Range("C24").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '10.46'
$ws.Range("E24").Value = '  +3.86%  '
$ws.Range("E25").Value = '  -1.76%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '161.79'
$ws.Range("E26").Value = '  +1.86%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '16.90'
$ws.Range("E27").Value = '  +6.25%  '
$ws.Range("E28").Value = '  +4.87%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.15'
$ws.Range("E29").Value = '  +5.54%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.993'
$ws.Range("E30").Value = '  -0.40%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.85'
$ws.Range("E31").Value = '  +10.82%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0516'
$ws.Range("E32").Value = '  +2.92%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.21'
$ws.Range("E33").Value = '  +6.27%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.60'
$ws.Range("E34").Value = '  +8.95%  '
$ws.Range("D35").Value = '1.575.22'
$ws.Range("E35").Value = '  +6.29%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.85'
$ws.Range("E36").Value = '  +5.96%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '88.75'
$ws.Range("E37").Value = '  +10.94%  '
$ws.Range("E38").Value = '  +2.63%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.628'
$ws.Range("E39").Value = '  +6.33%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0188'
$ws.Range("E40").Value = '  +4.82%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.84'
$ws.Range("E41").Value = '  +6.05%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.36'
$ws.Range("E42").Value = '  +2.77%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.922'
$ws.Range("E43").Value = '  +7.24%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.15'
$ws.Range("E44").Value = '  +6.16%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0521'
$ws.Range("E45").Value = '  +3.06%  '
$ws.Range("E46").Value = '  +2.99%  '
$ws.Range("D47").Value = '1.954.00'
$ws.Range("E47").Value = '  +7.77%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '54.13'
$ws.Range("E48").Value = '  +3.61%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '5.77'
$ws.Range("E49").Value = '  +5.83%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.997'
$ws.Range("E50").Value = '  -0.03%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '11.44'
$ws.Range("E51").Value = '  +22.11%  '
